$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D/E columns for rows with only price/volume changes
$ws.Range("D2").Value = '27.724.10'
$ws.Range("E2").Value = '  +0.36%  '

$ws.Range("D3").Value = '1.634.55'
$ws.Range("E3").Value = '  -0.94%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.15'
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.523'
$ws.Range("E6").Value = '  -2.30%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.28'
$ws.Range("E8").Value = '  -1.42%  '

$ws.Range("E9").Value = '  +1.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0610'
$ws.Range("E10").Value = '  -0.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("D12").Value = '1.866.60'
$ws.Range("E12").Value = '  -0.82%  '

$ws.Range("D13").Value = '1.626.72'
$ws.Range("E13").Value = '  -1.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.04'
$ws.Range("E14").Value = '  -0.23%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.560'
$ws.Range("E15").Value = '  -4.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.59'
$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("D17").Value = '27.707.95'
$ws.Range("E17").Value = '  +0.45%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.46'
$ws.Range("E18").Value = '  -1.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.66'
$ws.Range("E19").Value = '  +0.31%  '

$ws.Range("E20").Value = '  -0.50%  '

$ws.Range("E21").Value = '  +0.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.29'
$ws.Range("E22").Value = '  -1.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.99'
$ws.Range("E23").Value = '  +2.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.09'
$ws.Range("E24").Value = '  +5.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.27'
$ws.Range("E25").Value = '  +1.59%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.92'
$ws.Range("E26").Value = '  -2.09%  '

$ws.Range("E27").Value = '  -1.27%  '

$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.59'
$ws.Range("E29").Value = '  -0.32%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.18'
$ws.Range("E30").Value = '  -0.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0484'
$ws.Range("E31").Value = '  -0.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.30'
$ws.Range("E32").Value = '  -0.35%  '

$ws.Range("D33").Value = '1.466.69'
$ws.Range("E33").Value = '  +2.28%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.10'
$ws.Range("E34").Value = '  -3.11%  '

$ws.Range("E35").Value = '  -3.41%  '

$ws.Range("E36").Value = '  -0.62%  '

$ws.Range("E41").Value = '  +0.10%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '68.50'
$ws.Range("E42").Value = '  +5.14%  '

$ws.Range("E43").Value = '  -2.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.57'
$ws.Range("E44").Value = '  +1.34%  '

$ws.Range("E45").Value = '  +0.15%  '

$ws.Range("E46").Value = '  -1.14%  '

$ws.Range("D47").Value = '1.776.25'
$ws.Range("E47").Value = '  -0.87%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.69'
$ws.Range("E48").Value = '  -0.28%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.93'
$ws.Range("E49").Value = '  -1.28%  '

$ws.Range("D50").Value = '0.0₆0105'
$ws.Range("E50").Value = '  -1.29%  '

$ws.Range("E51").Value = '  -1.14%  '

# Rows 37-40 re-rank: TrustWalletToken moves to row 37, others shift down by one
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.940'
$ws.Range("E37").Value = '  +14.84%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.566'
$ws.Range("E38").Value = '  -1.42%  '

$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.876'
$ws.Range("E39").Value = '  -1.45%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0167'
$ws.Range("E40").Value = '  -0.65%  '
